$d = $word.ActiveDocument

# --- 1. Replace the "<<statementOfTruth.role>>" merge field with the
#        "<<{dateFormat(submittedOn, 'd MMMM yyyy', 'dd-MM-yyyy')}>>" field.
$lq = [char]0x2018
$rq = [char]0x2019
$newRoleText = "<<{dateFormat(submittedOn, " + $lq + "d MMMM yyyy" + $rq + ", " + $lq + "dd-MM-yyyy" + $rq + ")}>>"

$found = $d.Content.Find.Execute("<<statementOfTruth.role>>", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, $newRoleText, 2)
Write-Host "Replaced statementOfTruth.role: $found"

# --- 2. Update the paragraph formatting of the "<<statementOfTruth.name>>"
#        paragraph: its paragraph-mark run properties switch from shading
#        (white fill) to an explicit font colour (0A0A0A), matching the
#        surrounding paragraphs.
$found2 = $d.Content.Find.Execute("<<statementOfTruth.name>>", $true, $false, $false, $false, $false, `
                                   $true, 1, $false, "", 0)
Write-Host "Found statementOfTruth.name: $found2"
if ($found2) {
    $para = $d.Content.Find.Parent.Paragraphs(1)
    $pr = $d.Range($d.Content.Find.Parent.Start, $d.Content.Find.Parent.Start).Paragraphs(1)
}

$rng = $d.Content
$rng.Find.Execute("<<statementOfTruth.name>>") | Out-Null
$paraRange = $rng.Paragraphs(1).Range
$paraRange.ParagraphFormat.Shading.BackgroundPatternColor = -16777216
$paraRange.Font.Color = 657930
